$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 table: switch the applied table style from the custom
#    "Table_0" style to the built-in style {82294AF1-923D-423B-8A53-B66853A0918D}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{82294AF1-923D-423B-8A53-B66853A0918D}")

# ---------------------------------------------------------------------------
# 2) Theme colours: the deck's live theme (driving the slide master) was the
#    "Integral" palette; repoint it to the stock "Office Theme" palette by
#    rewriting each theme colour slot. RGB values are packed 0x00BBGGRR, as
#    PowerPoint's object model expects.
#      1 dk1      000000
#      2 lt1      FFFFFF
#      3 dk2      44546A
#      4 lt2      E7E6E6
#      5 accent1  5B9BD5
#      6 accent2  ED7D31
#      7 accent3  A5A5A5
#      8 accent4  FFC000
#      9 accent5  4472C4
#     10 accent6  70AD47
#     11 hlink    0563C1
#     12 folHlink 954F72
# ---------------------------------------------------------------------------
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB = 0
$colorScheme.Item(2).RGB = 16777215
$colorScheme.Item(3).RGB = 6968388
$colorScheme.Item(4).RGB = 15132391
$colorScheme.Item(5).RGB = 13998939
$colorScheme.Item(6).RGB = 3243501
$colorScheme.Item(7).RGB = 10855845
$colorScheme.Item(8).RGB = 49407
$colorScheme.Item(9).RGB = 12874308
$colorScheme.Item(10).RGB = 4697456
$colorScheme.Item(11).RGB = 12673797
$colorScheme.Item(12).RGB = 7491477
